$wb = $excel.ActiveWorkbook

# --- Sheet "traceback results" ---
$ws1 = $wb.Worksheets.Item("traceback results")

# Row 2 (scenario 7, Chain 1_5_1 -> Chain 1)
$ws1.Range("A2").Value = 7
$ws1.Range("D2").Value = 0.9007217937885342
$ws1.Range("E2").Value = 0.06886383948735718
$ws1.Range("F2").Value = -17.60655854649039
$ws1.Range("G2").Value = 17.43300656077044
$ws1.Range("H2").Value = 0.3471039714398998
$ws1.Range("I2").Value = 0.8406734470115096

# Row 3 (scenario 7, Chain 1_5_1 -> Chain 2)
$ws1.Range("A3").Value = 7
$ws1.Range("B3").Value = "Chain 1_5_1"
$ws1.Range("C3").Value = "Chain 2"
$ws1.Range("D3").Value = 0.001000000000038916
$ws1.Range("E3").Value = 49.99999986057535
$ws1.Range("F3").Value = -17.60655854649039
$ws1.Range("G3").Value = 17.60658123503512
$ws1.Range("H3").Value = [double]"-4.537708944951646e-05"
$ws1.Range("I3").Value = 1
$ws1.Range("J3").Value = "Fail to reject the null hypothesis."

# Row 4 (scenario 7, Chain 2_5_1 -> Chain 1)
$ws1.Range("A4").Value = 7
$ws1.Range("B4").Value = "Chain 2_5_1"
$ws1.Range("C4").Value = "Chain 1"
$ws1.Range("D4").Value = 0.001000000000003481
$ws1.Range("E4").Value = 49.99999784264019
$ws1.Range("F4").Value = -20.49413135993761
$ws1.Range("G4").Value = 20.49415386146449
$ws1.Range("H4").Value = [double]"-4.500305375643165e-05"
$ws1.Range("I4").Value = 1
$ws1.Range("J4").Value = "Fail to reject the null hypothesis."

# Row 5 (scenario 7, Chain 2_5_1 -> Chain 2)
$ws1.Range("A5").Value = 7
$ws1.Range("D5").Value = 0.8486166910306896
$ws1.Range("E5").Value = 0.1657004664173341
$ws1.Range("F5").Value = -20.49413135993761
$ws1.Range("G5").Value = 20.33119640499809
$ws1.Range("H5").Value = 0.3258699098790458
$ws1.Range("I5").Value = 0.8496464519221111

# --- Sheet "flow results" ---
$ws2 = $wb.Worksheets.Item("flow results")
$ws2.Range("A2").Value = 7
$ws2.Range("B2").Value = 16.83176100781197
$ws2.Range("C2").Value = 0.0001876266873183097
